{"js": "// Updated figures for presentations:\n// Every table cell shaded dark gray (A9A9A9) becomes white (FFFFFF).\n\nconst GRAY = \"#A9A9A9\";\nconst WHITE = \"#FFFFFF\";\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Load the rows collection for every table in the document.\nfor (const table of tables.items) {\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\n// Load the cells collection for every row of every table.\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    row.cells.load(\"items\");\n  }\n}\nawait context.sync();\n\n// Load the current shading color of every cell so we can test it.\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    for (const cell of row.cells.items) {\n      cell.load(\"shadingColor\");\n    }\n  }\n}\nawait context.sync();\n\n// Flip every gray-shaded cell to white; leave every other cell untouched.\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    for (const cell of row.cells.items) {\n      const current = (cell.shadingColor || \"\").toUpperCase();\n      if (current === GRAY) {\n        cell.shadingColor = WHITE;\n      }\n    }\n  }\n}\nawait context.sync();\n", "ps1": "# Updated figures for presentations:\n# Every table cell shaded dark gray (A9A9A9) becomes white (FFFFFF).\n\n$d = $word.ActiveDocument\n\n$GRAY  = 0xA9A9A9\n$WHITE = 0xFFFFFF\n\nforeach ($table in $d.Tables) {\n    foreach ($row in $table.Rows) {\n        foreach ($cell in $row.Cells) {\n            if ($cell.Shading.BackgroundPatternColor -eq $GRAY) {\n                $cell.Shading.BackgroundPatternColor = $WHITE\n            }\n        }\n    }\n}\n"}
